# This script reproduces the symbol-list refresh committed by the
# "Updated symbol list" GitHub Actions job: the Price (column D) and
# Volume(1h) (column E) cells on the active sheet are refreshed with new
# quoted values. The source data is free-form text (prices have varying
# decimal precision, e.g. "5.064" or "0.00000000750", and volumes are
# percentage strings like "1.09%"), so every target cell is forced to the
# Text number format before the value is written. This guarantees Excel
# stores the exact characters instead of silently re-interpreting the
# text as a number/percentage and rounding or reformatting it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue "D2" "304.55"
Set-TextValue "E2" "1.09%"
Set-TextValue "D3" "35.66"
Set-TextValue "E3" "1.42%"
Set-TextValue "D4" "5.064"
Set-TextValue "E4" "0.46%"
Set-TextValue "E5" "0.83%"
Set-TextValue "D6" "1.915"
Set-TextValue "E6" "0.56%"
Set-TextValue "D7" "4.178"
Set-TextValue "E7" "3.08%"
Set-TextValue "D8" "7.736"
Set-TextValue "E8" "-0.83%"
Set-TextValue "D9" "0.9284"
Set-TextValue "E9" "0.70%"
Set-TextValue "D10" "0.1391"
Set-TextValue "E10" "8.67%"
Set-TextValue "D11" "0.1897"
Set-TextValue "E11" "2.88%"
Set-TextValue "D12" "0.09235"
Set-TextValue "E12" "-6.35%"
Set-TextValue "D13" "0.03592"
Set-TextValue "E13" "0.43%"
Set-TextValue "D14" "0.09813"
Set-TextValue "E14" "-0.37%"
Set-TextValue "D15" "0.001408"
Set-TextValue "E15" "1.58%"
Set-TextValue "D16" "0.005903"
Set-TextValue "E16" "2.07%"
Set-TextValue "D17" "3.553"
Set-TextValue "E17" "1.39%"
Set-TextValue "D18" "3.008"
Set-TextValue "E18" "3.22%"
Set-TextValue "D19" "0.3469"
Set-TextValue "E19" "2.05%"
Set-TextValue "E20" "0.14%"
Set-TextValue "D21" "4.902"
Set-TextValue "E21" "-3.01%"
Set-TextValue "E22" "4.42%"
Set-TextValue "E23" "-1.41%"
Set-TextValue "E24" "0.55%"
Set-TextValue "D25" "0.004783"
Set-TextValue "E25" "0.05%"
Set-TextValue "E26" "24.76%"
Set-TextValue "D27" "0.0003130"
Set-TextValue "E27" "4.22%"
Set-TextValue "E39" "5.00%"
Set-TextValue "D40" "0.04895"
Set-TextValue "E40" "4.28%"
Set-TextValue "D41" "0.007643"
Set-TextValue "E41" "1.89%"
Set-TextValue "D42" "0.009402"
Set-TextValue "E42" "-8.85%"
Set-TextValue "D43" "0.1372"
Set-TextValue "E43" "3.82%"
Set-TextValue "E44" "-0.55%"
Set-TextValue "E45" "6.13%"
Set-TextValue "D46" "0.00006370"
Set-TextValue "E46" "1.95%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.06%"
Set-TextValue "D48" "63.57"
Set-TextValue "E48" "-1.41%"
Set-TextValue "D49" "0.001191"
Set-TextValue "E49" "-20.02%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "-0.06%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "-0.06%"
